# Update "想去人数" (interest/attendance counter) figures on the
# "展览" (Exhibition) and "全部类型" (All types) sheets to match the
# newly scraped counts (output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F7").Value = 2668
$wsExhibit.Range("F8").Value = 1167
$wsExhibit.Range("F11").Value = 9915
$wsExhibit.Range("F13").Value = 255
$wsExhibit.Range("F14").Value = 4
$wsExhibit.Range("F15").Value = 612
$wsExhibit.Range("F16").Value = 11732
$wsExhibit.Range("F17").Value = 12064

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 2668
$wsAll.Range("F9").Value = 1167
$wsAll.Range("F12").Value = 9916
$wsAll.Range("F14").Value = 255
$wsAll.Range("F15").Value = 4
$wsAll.Range("F16").Value = 612
$wsAll.Range("F17").Value = 11732
$wsAll.Range("F18").Value = 12064
